$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.231.53"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "1.905.26"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.33"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5356"
$ws.Range("E7").Value = "  +2.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3820"
$ws.Range("E8").Value = "  +1.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07307"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("E10").Value = "  +5.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9059"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08205"
$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.74"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.355"
$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.004"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.87"
$ws.Range("E16").Value = "  +2.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008670"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D19").Value = "27.256.97"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("D21").Value = "1.071.41"
$ws.Range("E21").Value = "  -43.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  +0.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.519"
$ws.Range("E23").Value = "  +1.81%  "

$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.297"
$ws.Range("E25").Value = "  +0.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.39"
$ws.Range("E26").Value = "  +0.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.747"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.81"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.820"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.721"
$ws.Range("E30").Value = "  -4.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09225"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8294"
$ws.Range("E32").Value = "  +4.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05081"
$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.218"
$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.997"
$ws.Range("E35").Value = "  +1.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.325"
$ws.Range("E36").Value = "  -3.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("E37").Value = "  +3.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5906"
$ws.Range("E38").Value = "  +4.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02002"
$ws.Range("E39").Value = "  +0.80%  "

$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.364"
$ws.Range("E41").Value = "  +4.86%  "

$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.24"
$ws.Range("E43").Value = "  +1.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5090"
$ws.Range("E44").Value = "  +4.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1527"
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.12"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.642"
$ws.Range("E48").Value = "  +1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.33"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  +3.63%  "

$ws.Range("E51").Value = "  +0.25%  "
